$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that directly follows the
#    "Play Chilli Master slot for free" heading paragraph at the top of the
#    document. Deleting the whole paragraph Range (including its end-of-
#    paragraph mark) removes the paragraph entirely.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Near the end of the document, the final paragraph (an italic paragraph
#    that used to hold the feature-image prompt text) is replaced by two
#    paragraphs:
#      - a new bold paragraph containing "Play Chilli Master slot for free"
#      - the existing italic paragraph, but with its text swapped for the
#        "Meta description" copy that used to live at the top of the doc.
#    InsertXML lets us replace the whole paragraph (including its paragraph
#    mark) with exactly the run/paragraph structure we need.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Chilli Master slot for free</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Chilli Master slot, a Mexican-themed game with free spins and potential for big wins. Play this game for free today.</w:t></w:r></w:p>'

$lastPara.Range.InsertXML($xml) | Out-Null
